$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the discontinued product row ("DIFLUSTERO 0.05% EYE EMULSION 5 ML").
# Deleting the entire row shifts every row below it up by one, which keeps
# the running item numbers (column A) sequential and moves the totals /
# footer rows up to their new positions automatically.
$ws.Rows.Item(25).Delete()

# The grand-total cell is a literal value (not a formula), so it must be
# corrected by hand to remove the deleted item's price (45.00).
$total = $ws.Cells.Item(73, 16)
$total.Value2 = $total.Value2 - 45

# Update the generated-on timestamp in the report footer.
$ws.Cells.Item(74, 1).Value2 = "Monday, 28 July, 2025 7:35 PM"
